# Auto-generated Excel COM-interop edit script
# Updates currentAveragePrice / Leve profit calculations across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 885.5
$ws.Range("I28").Value = 1214.091
$ws.Range("J28").Value = 162.6
$ws.Range("K28").Value = 1214.091
$ws.Range("L28").Value = 162.6
$ws.Range("M28").Value = -729.0909999999999
$ws.Range("N28").Value = -1132.6

$ws.Range("H34").Value = 29500
$ws.Range("J34").Value = 29500
$ws.Range("L34").Value = 29500
$ws.Range("N34").Value = -29906

$ws.Range("H36").Value = 29500
$ws.Range("J36").Value = 29500
$ws.Range("L36").Value = 29500
$ws.Range("N36").Value = -30930

$ws.Range("H107").Value = 301
$ws.Range("I107").Value = 212.90909
$ws.Range("J107").Value = 494.8
$ws.Range("K107").Value = 212.90909
$ws.Range("L107").Value = 494.8
$ws.Range("M107").Value = 1707.09091
$ws.Range("N107").Value = -4334.8

$ws.Range("H133").Value = 45411.766
$ws.Range("J133").Value = 45411.766
$ws.Range("L133").Value = 45411.766
$ws.Range("N133").Value = -55531.766

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1083.4445
$ws.Range("I97").Value = 1210
$ws.Range("J97").Value = 982.2
$ws.Range("K97").Value = 1210
$ws.Range("L97").Value = 982.2
$ws.Range("M97").Value = -714
$ws.Range("N97").Value = -1974.2

$ws.Range("H132").Value = 2110.4595
$ws.Range("I132").Value = 1236.9048
$ws.Range("J132").Value = 3257
$ws.Range("K132").Value = 3710.7144
$ws.Range("L132").Value = 9771
$ws.Range("M132").Value = -1180.7144
$ws.Range("N132").Value = -14831

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4037.318
$ws.Range("I105").Value = 5071.5386
$ws.Range("K105").Value = 5071.5386
$ws.Range("M105").Value = -3324.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 114.117645
$ws.Range("I7").Value = 125.833336
$ws.Range("J7").Value = 86
$ws.Range("K7").Value = 125.833336
$ws.Range("L7").Value = 86
$ws.Range("M7").Value = -12.833336
$ws.Range("N7").Value = -312

$ws.Range("H22").Value = 338.84616
$ws.Range("I22").Value = 375.875
$ws.Range("J22").Value = 279.6
$ws.Range("K22").Value = 375.875
$ws.Range("L22").Value = 279.6
$ws.Range("M22").Value = -25.875
$ws.Range("N22").Value = -979.6

$ws.Range("H99").Value = 2636.3447
$ws.Range("I99").Value = 2275.5557
$ws.Range("J99").Value = 3226.7273
$ws.Range("K99").Value = 2275.5557
$ws.Range("L99").Value = 3226.7273
$ws.Range("M99").Value = -777.5556999999999
$ws.Range("N99").Value = -6222.7273

$ws.Range("H126").Value = 2636.3447
$ws.Range("I126").Value = 2275.5557
$ws.Range("J126").Value = 3226.7273
$ws.Range("K126").Value = 6826.6671
$ws.Range("L126").Value = 9680.1819
$ws.Range("M126").Value = -4356.6671
$ws.Range("N126").Value = -14620.1819

$ws.Range("H140").Value = 59864
$ws.Range("J140").Value = 59864
$ws.Range("L140").Value = 59864
$ws.Range("N140").Value = -70224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 724.55554
$ws.Range("J113").Value = 586.8333
$ws.Range("L113").Value = 1760.4999
$ws.Range("N113").Value = -6100.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1838.2667
$ws.Range("I97").Value = 1839.4445
$ws.Range("J97").Value = 1836.5
$ws.Range("K97").Value = 1839.4445
$ws.Range("L97").Value = 1836.5
$ws.Range("M97").Value = -1343.4445
$ws.Range("N97").Value = -2828.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 506414.6
$ws.Range("I40").Value = 595517.2
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 595517.2
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -595381.2
$ws.Range("N40").Value = -1772

$ws.Range("H100").Value = 4832283
$ws.Range("I100").Value = 6537230
$ws.Range("J100").Value = 1599
$ws.Range("K100").Value = 6537230
$ws.Range("L100").Value = 1599
$ws.Range("M100").Value = -6536689
$ws.Range("N100").Value = -2681

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 42000
$ws.Range("J92").Value = 42000
$ws.Range("L92").Value = 42000
$ws.Range("N92").Value = -46992

$ws.Range("H93").Value = 25311
$ws.Range("J93").Value = 25311
$ws.Range("L93").Value = 25311
$ws.Range("N93").Value = -30303

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H96").Value = 3146.889
$ws.Range("I96").Value = 2467.9092
$ws.Range("K96").Value = 2467.9092
$ws.Range("M96").Value = -1094.9092

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H100").Value = 385.07144
$ws.Range("I100").Value = 321.22223
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 642.44446
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -101.44446
$ws.Range("N100").Value = -2082

$ws.Range("H101").Value = 36602
$ws.Range("J101").Value = 36602
$ws.Range("L101").Value = 36602
$ws.Range("N101").Value = -43092

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H104").Value = 19900
$ws.Range("J104").Value = 19900
$ws.Range("L104").Value = 19900
$ws.Range("N104").Value = -26888

$ws.Range("H105").Value = 42326.625
$ws.Range("J105").Value = 42326.625
$ws.Range("L105").Value = 42326.625
$ws.Range("N105").Value = -49314.625

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H122").Value = 1138.4166
$ws.Range("I122").Value = 926
$ws.Range("J122").Value = 1350.8334
$ws.Range("K122").Value = 2778
$ws.Range("L122").Value = 4052.5002
$ws.Range("M122").Value = -328
$ws.Range("N122").Value = -8952.5002
